$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 29
$ws1.Range("F5").Value = 15999
$ws1.Range("F9").Value = 15516
$ws1.Range("F10").Value = 63
$ws1.Range("F11").Value = 9140
$ws1.Range("F20").Value = 71
$ws1.Range("F29").Value = 506
$ws1.Range("F34").Value = 59
$ws1.Range("F36").Value = 336
$ws1.Range("F39").Value = 5620

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 2

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 29
$ws4.Range("F5").Value = 15999
$ws4.Range("F9").Value = 15516
$ws4.Range("F10").Value = 63
$ws4.Range("F11").Value = 9140
$ws4.Range("F20").Value = 71
$ws4.Range("F29").Value = 506
$ws4.Range("F36").Value = 59
$ws4.Range("F38").Value = 336
$ws4.Range("F41").Value = 5620
$ws4.Range("F42").Value = 2

$wb.Save()
